# Ingreso: cabecera compra, masivo con compra y proveedores auto-upsert
#
# - Adds purchase-header columns (Tipo doc, N° documento, Fecha doc,
#   RUT proveedor, Proveedor) to the "Ingreso EPP" sheet.
# - Sets explicit column widths for the full A:K header range.
# - Freezes the header row (row 1).
# - Drops the old per-cell validation comments and their data validations
#   (now presumably handled elsewhere / superseded by the new workflow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells (G1:K1) -------------------------------------------
$ws.Range("G1").Value = "Tipo doc"
$ws.Range("H1").Value = "N° documento"
$ws.Range("I1").Value = "Fecha doc"
$ws.Range("J1").Value = "RUT proveedor"
$ws.Range("K1").Value = "Proveedor"

# --- 2. Column widths (A:K) -------------------------------------------------
# ColumnWidth is in "characters"; the stored XML width adds ~5/6 of a
# character of padding, so back it out to land on the exact target widths.
function Set-ColWidth($col, $target) {
    $ws.Columns.Item($col).ColumnWidth = $target - (5/6)
}

Set-ColWidth 1 18
Set-ColWidth 2 28
Set-ColWidth 3 14
Set-ColWidth 4 10
Set-ColWidth 5 14
Set-ColWidth 6 12
Set-ColWidth 7 12
Set-ColWidth 8 14
Set-ColWidth 9 12
Set-ColWidth 10 14
Set-ColWidth 11 26

# --- 3. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# --- 4. Remove old data validations ----------------------------------------
$ws.Range("A2:A101").Validation.Delete()
$ws.Range("D2:D101").Validation.Delete()
$ws.Range("E2:E101").Validation.Delete()
$ws.Range("F2:F101").Validation.Delete()

# --- 5. Remove the legacy cell comments (and their VML/legacyDrawing) ------
$ws.Range("A1").Comment.Delete()
$ws.Range("D1").Comment.Delete()
$ws.Range("E1").Comment.Delete()
$ws.Range("F1").Comment.Delete()
